$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row cells: "_old" -> "_FV2310", "_new" -> "_FV2404"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $txt = $cell.Value2
    if ($txt -match "_old$") {
        $cell.Value = ($txt -replace "_old$", "_FV2310")
    } elseif ($txt -match "_new$") {
        $cell.Value = ($txt -replace "_new$", "_FV2404")
    }
}

# 2) Turn the data range into an Excel Table ("Table1") with headers
$range = $ws.Range("A1:U80")
$lo = $ws.ListObjects.Add(1, $range, $null, 1)
$lo.Name = "Table1"

# 3) Freeze the header row (split after row 1, keep top-left at A2)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
